# Update "想去人数" (number of people interested) figures across the four
# sheets of the workbook (展览, 演出, 本地生活, 全部类型). Only column F
# values change; everything else stays the same.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1192
$ws1.Range("F3").Value  = 2011
$ws1.Range("F7").Value  = 57
$ws1.Range("F8").Value  = 145
$ws1.Range("F9").Value  = 358
$ws1.Range("F10").Value = 142
$ws1.Range("F12").Value = 890
$ws1.Range("F13").Value = 280
$ws1.Range("F14").Value = 145
$ws1.Range("F17").Value = 349
$ws1.Range("F18").Value = 282
$ws1.Range("F19").Value = 719
$ws1.Range("F20").Value = 93
$ws1.Range("F21").Value = 679
$ws1.Range("F22").Value = 217
$ws1.Range("F23").Value = 51
$ws1.Range("F24").Value = 931
$ws1.Range("F25").Value = 386
$ws1.Range("F26").Value = 207
$ws1.Range("F27").Value = 63
$ws1.Range("F28").Value = 320
$ws1.Range("F30").Value = 25
$ws1.Range("F31").Value = 437

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 341
$ws2.Range("F6").Value  = 34
$ws2.Range("F9").Value  = 6
$ws2.Range("F11").Value = 136

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 337

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 337
$ws4.Range("F3").Value  = 1192
$ws4.Range("F4").Value  = 2011
$ws4.Range("F9").Value  = 57
$ws4.Range("F10").Value = 145
$ws4.Range("F11").Value = 358
$ws4.Range("F12").Value = 142
$ws4.Range("F14").Value = 890
$ws4.Range("F15").Value = 280
$ws4.Range("F16").Value = 145
$ws4.Range("F19").Value = 341
$ws4.Range("F22").Value = 349
$ws4.Range("F23").Value = 34
$ws4.Range("F25").Value = 282
$ws4.Range("F26").Value = 719
$ws4.Range("F27").Value = 93
$ws4.Range("F28").Value = 679
$ws4.Range("F29").Value = 217
$ws4.Range("F30").Value = 51
$ws4.Range("F31").Value = 931
$ws4.Range("F32").Value = 386
$ws4.Range("F34").Value = 6
$ws4.Range("F35").Value = 207
$ws4.Range("F36").Value = 63
$ws4.Range("F37").Value = 320
$ws4.Range("F39").Value = 136
$ws4.Range("F41").Value = 25
$ws4.Range("F43").Value = 437
